$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2, 0.0082523160758350891, -0.25397548895418398, -0.22162692975249562),
    @(3, 0.0027851325529865451, 0.1067070414303771, 0.11762458560184386),
    @(4, 0.0092740682930517633, 0.18214944453524412, 0.21850320806797785),
    @(5, 0.0040202977998172783, -0.09746380117166234, -0.081704486250136987),
    @(6, 0.0023298523953054141, 0.064145004910871065, 0.073277879997875064),
    @(7, 0.0051691716355643288, 0.058109867463255749, 0.078372695677617621),
    @(8, 0.0059580004339406063, -0.28112036946152708, -0.25776538305651242),
    @(9, 0.0025454379104683776, 0.21350238818572695, 0.22348034445695178),
    @(10, 0.0088900658953239391, 0.15251024076294284, 0.18735873908499667),
    @(11, 0.0037871233486962493, -0.12494901436645819, -0.11010372865118645),
    @(12, 0.00280768030287881, 0.11550867635245136, 0.12651460683203969),
    @(13, 0.0063040747841167106, 0.065518324751005455, 0.090229902041687141),
    @(14, 0.0047049622321926965, -0.25341537793449587, -0.23497222235112325),
    @(15, 0.0029222914826086892, 0.27698216655579316, 0.2884373650917012),
    @(16, 0.0083136041074445399, 0.069615900194180513, 0.10220470461922851),
    @(17, 0.00323874158855541, -0.12813617583860956, -0.11544051218757213),
    @(18, 0.0027580862441502567, 0.1611856754927434, 0.17199720037636432),
    @(19, 0.0061326191804874589, 0.041723728310545723, 0.065763210401519459),
    @(20, 0.0040277527922005871, -0.21276732343374988, -0.19697878619754355),
    @(21, 0.0036250026856564557, 0.29845119496276112, 0.31266097715065083),
    @(22, 0.0072066990102844612, 0.0075693725891402813, 0.035819178757567342),
    @(23, 0.003161788922094478, -0.12054785092354448, -0.10815383689280733),
    @(24, 0.0031129638886849284, 0.19316793391698828, 0.20537055688271777),
    @(25, 0.0055689767953412996, 0.031873528761391906, 0.05370356809639823),
    @(26, 0.003534003659928607, -0.17808173808908814, -0.1642286663499988),
    @(27, 0.003247886759319064, 0.30016056904421362, 0.31289208055549322),
    @(28, 0.0079306676462223419, -0.044302961229235521, -0.013215243610908953),
    @(29, 0.0032098531826458064, -0.11571907708939545, -0.10313665417548228),
    @(30, 0.0031593090285371946, 0.20438572368626812, 0.21677001668998513),
    @(31, 0.0071093801118957307, 0.019026200697609987, 0.046894524304206014),
    @(32, 0.0035922697164041911, -0.14521359203014819, -0.13113212101986896),
    @(33, 0.0034757582953194445, 0.26642151810686943, 0.28004627168558194),
    @(34, 0.0082906691727075623, -0.081179293849689318, -0.048680392924132634),
    @(35, 0.0029269237933898049, -0.10399238102469097, -0.092519023550169555),
    @(36, 0.0035141449846811507, 0.2026636966650765, 0.21643892433502826),
    @(37, 0.0063272972513269219, 0.018025218766694193, 0.04282782667058993),
    @(38, 0.0026416535607332692, -0.11336707308547141, -0.10301195752585761),
    @(39, 0.002884914085758548, 0.18141649569169188, 0.19272517718635956),
    @(40, 0.0067461448544357773, -0.088293204270811929, -0.061848741382882891),
    @(41, 0.0032249015768308466, -0.095907155611213143, -0.083265743937055367),
    @(42, 0.0036030286525759506, 0.1920106744319782, 0.20613432049864625),
    @(43, 0.0076878773670840421, 0.0058942125758719086, 0.036030209096192131),
    @(44, 0.0035418169100142065, -0.10659149224718201, -0.092707793059915777),
    @(45, 0.0030637085376993426, 0.07536829197117112, 0.087377836455136637),
    @(46, 0.007146571937682992, -0.073617871153966938, -0.045603759322717333),
    @(47, 0.0033311083971766238, -0.084250985501828468, -0.071193249761150296),
    @(48, 0.0036975602429509752, 0.1566830423771452, 0.17117724634199194),
    @(49, 0.0084122778742870222, -0.010407613394357505, 0.022567987625620258),
    @(50, 0.003757050738416205, -0.10204444197592427, -0.087317039738958324),
    @(51, 0.0038441029296961706, -0.019268243536715867, -0.0041996021933731554),
    @(52, 0.0072747871982678071, -0.014405053911784234, 0.014111653664644468),
    @(53, 0.0039416519128551162, -0.070837480805230266, -0.055386452822044588),
    @(54, 0.0050334218341020545, 0.096348221565284622, 0.11607891908229474),
    @(55, 0.010256166945429537, -0.030462805076227922, 0.0097407253160004415)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 3).Value = $row[1]
    $ws.Cells.Item($r, 4).Value = $row[2]
    $ws.Cells.Item($r, 5).Value = $row[3]
}
